$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3, 4 and 5 (the three data records below the header row) need to be
# cyclically rotated:
#   new row 3 <- old row 4
#   new row 4 <- old row 5
#   new row 5 <- old row 3
# Capture the full row contents (columns A:AY) first so the rotation does
# not clobber data before it has been read.

$row3 = $ws.Range("A3:AY3").Value()
$row4 = $ws.Range("A4:AY4").Value()
$row5 = $ws.Range("A5:AY5").Value()

$ws.Range("A3:AY3").Value = $row4
$ws.Range("A4:AY4").Value = $row5
$ws.Range("A5:AY5").Value = $row3

# Columns Y and AA hold plain text dates ("2020-02-29"), and column I
# sometimes holds a digit stored as text ("1"). A bulk Value assignment lets
# Excel's automatic type detection reinterpret these as a real date/number,
# which is not what the source data contained. Re-apply them as text
# (leading apostrophe forces text, and is not retained in the stored value).
$ws.Range("Y3").Value = "'" + $row4.GetValue(1, 25)
$ws.Range("AA3").Value = "'" + $row4.GetValue(1, 27)

$ws.Range("Y4").Value = "'" + $row5.GetValue(1, 25)
$ws.Range("AA4").Value = "'" + $row5.GetValue(1, 27)

$ws.Range("Y5").Value = "'" + $row3.GetValue(1, 25)
$ws.Range("AA5").Value = "'" + $row3.GetValue(1, 27)
if ($row3.GetValue(1, 9) -ne $null) {
    $ws.Range("I5").Value = "'" + $row3.GetValue(1, 9)
}
